$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> (Price, Volume) updates as seen in the diff.
# $null for D means the Price column is unchanged for that row.
$updates = @(
    @{ Row = 2;  D = "29.656.65"; E = "  +0.53%  " },
    @{ Row = 3;  D = "1.612.11";  E = "  +0.59%  " },
    @{ Row = 4;  D = $null;       E = "  -0.57%  " },
    @{ Row = 5;  D = "212.56";    E = "  -0.13%  " },
    @{ Row = 6;  D = $null;       E = "  -0.76%  " },
    @{ Row = 7;  D = $null;       E = "  -0.56%  " },
    @{ Row = 8;  D = $null;       E = "  +7.70%  " },
    @{ Row = 9;  D = $null;       E = "  +2.85%  " },
    @{ Row = 10; D = $null;       E = "  +1.61%  " },
    @{ Row = 11; D = $null;       E = "  -0.46%  " },
    @{ Row = 12; D = "1.842.75";  E = "  +0.55%  " },
    @{ Row = 13; D = "1.609.24";  E = "  +0.51%  " },
    @{ Row = 14; D = $null;       E = "  +6.21%  " },
    @{ Row = 15; D = $null;       E = "  +3.37%  " },
    @{ Row = 16; D = "29.663.01"; E = "  +0.43%  " },
    @{ Row = 17; D = "8.85";      E = "  +16.24%  " },
    @{ Row = 18; D = "64.52";     E = "  +1.83%  " },
    @{ Row = 19; D = "240.72";    E = "  -0.69%  " },
    @{ Row = 20; D = $null;       E = "  +1.91%  " },
    @{ Row = 21; D = $null;       E = "  -0.39%  " },
    @{ Row = 22; D = $null;       E = "  +2.44%  " },
    @{ Row = 23; D = "9.63";      E = "  +5.05%  " },
    @{ Row = 24; D = $null;       E = "  +1.03%  " },
    @{ Row = 25; D = "156.45";    E = "  +1.41%  " },
    @{ Row = 26; D = "15.60";     E = "  +2.12%  " },
    @{ Row = 27; D = $null;       E = "  +0.81%  " },
    @{ Row = 28; D = $null;       E = "  +2.46%  " },
    @{ Row = 29; D = "0.994";     E = "  -0.52%  " },
    @{ Row = 30; D = $null;       E = "  +1.86%  " },
    @{ Row = 31; D = "3.29";      E = "  +2.23%  " },
    @{ Row = 32; D = $null;       E = "  +0.57%  " },
    @{ Row = 33; D = $null;       E = "  +2.81%  " },
    @{ Row = 34; D = "1.437.72";  E = "  +1.58%  " },
    @{ Row = 35; D = "1.61";      E = "  +5.96%  " },
    @{ Row = 36; D = $null;       E = "  +2.14%  " },
    @{ Row = 37; D = $null;       E = "  +3.75%  " },
    @{ Row = 38; D = "2.28";      E = "  -0.92%  " },
    @{ Row = 39; D = $null;       E = "  +3.22%  " },
    @{ Row = 40; D = "0.553";     E = "  +3.55%  " },
    @{ Row = 41; D = $null;       E = "  +5.91%  " },
    @{ Row = 42; D = $null;       E = "  +1.80%  " },
    @{ Row = 43; D = $null;       E = "  +3.99%  " },
    @{ Row = 44; D = "53.92";     E = "  +2.34%  " },
    @{ Row = 45; D = "69.66";     E = "  +6.23%  " },
    @{ Row = 46; D = $null;       E = "  -0.48%  " },
    @{ Row = 47; D = $null;       E = "  +20.33%  " },
    @{ Row = 48; D = $null;       E = "  +2.96%  " },
    @{ Row = 49; D = "1.751.57";  E = "  +0.44%  " },
    @{ Row = 50; D = "87.51";     E = "  +1.36%  " },
    @{ Row = 51; D = $null;       E = "  -1.48%  " }
)

# First, force the Price column cells to text format so that purely
# numeric-looking values (e.g. "212.56") are kept as text, matching the
# original inline-string cell type instead of being auto-converted to a
# number by Excel.
foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($u.Row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
